$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J - copy formatting (bold, border, centered)
# from an existing header cell (H1) so the new headers share the same style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Data rows 2-68 for columns I (9) and J (10)
$data = @"
2,4,6
3,7,7
4,8,8
5,7,7
6,7,7
7,7,7
8,8,8
9,8,8
10,7,7
11,10,10
12,8,8
13,6,6
14,7,7
15,8,8
16,8,8
17,5,6
18,8,8
19,7,7
20,7,7
21,9,9
22,7,7
23,7,7
24,7,7
25,9,9
26,7,7
27,6,6
28,10,10
29,7,7
30,6,6
31,9,9
32,6,6
33,9,9
34,9,9
35,7,7
36,6,6
37,6,7
38,8,8
39,7,7
40,6,6
41,6,7
42,6,8
43,9,9
44,6,6
45,7,7
46,9,9
47,9,9
48,7,7
49,6,6
50,9,9
51,9,9
52,6,6
53,9,9
54,8,8
55,5,6
56,8,8
57,9,9
58,8,8
59,9,9
60,9,9
61,5,5
62,6,6
63,6,6
64,8,8
65,7,7
66,5,5
67,3,4
68,4,4
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split ","
    $row = [int]$parts[0]
    $iVal = [int]$parts[1]
    $jVal = [int]$parts[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
